$d = $word.ActiveDocument

# --- Step 1: "company" -> "airline" ------------------------------------
$r = $d.Content
$r.Find.Execute("company")
$companyStart = $r.Start
$r.Text = "airline"

# --- Step 2: move " flight" from after the _GoBack bookmark to before it,
#             so "...seats even in" + " flight" => "...seats even in flight"
#             and the old " flight. " becomes just ". "
$moveRange = $d.Content
$moveRange.Find.Execute(" flight. ")
$flightStart = $moveRange.Start

$insertPoint = $d.Range($flightStart, $flightStart)
$insertPoint.InsertBefore(" flight")

$oldFlightText = $d.Range($flightStart + 7, $flightStart + 14)
$oldFlightText.Text = ""

# --- Step 3: force run boundaries so "airline" and the trailing clause
#             each become their own <w:r>, instead of being re-merged with
#             their neighbours. Toggling a character property on and back
#             off is enough to break the run at that point.
$airlineRange = $d.Range($companyStart, $companyStart + 7)
$airlineRange.Bold = 1
$airlineRange.Bold = 0

$tailStart = $companyStart + 7
$full = $d.Content.Text
$tailEnd = $full.IndexOf(" flight. Unlike") + 1
# the bookmark sits between "even in" and "flight" now, but Range
# positions are not affected by zero-width bookmarks, so locate the tail
# end via Find instead of raw text offsets.
$tailRange = $d.Content
$tailRange.Find.Execute(", but they can negotiate the seats even in flight")
Write-Output "tail: $($tailRange.Start) - $($tailRange.End) [$($tailRange.Text)]"
$tailRange.Bold = 1
$tailRange.Bold = 0
